$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The data rows for Níspero (rows 2-5) got reordered/updated during the
# weekly refresh: row 2 now holds what used to be row 4's record, row 3
# now holds what used to be row 5's record, and rows 4/5 now hold the
# former rows 2/3 records. Only columns D, K, L, M, N, O, P, Q, S change;
# everything else stays the same.

$ws.Range("D2").Value = 44505
$ws.Range("K2").Value = "Californiana(o)"
$ws.Range("L2").Value = "Primera"
$ws.Range("M2").Value = 100
$ws.Range("N2").Value = 15000
$ws.Range("O2").Value = 15000
$ws.Range("P2").Value = 15000
$ws.Range("Q2").Value = "$/bandeja 10 kilos"
$ws.Range("S2").Value = 1500

$ws.Range("D3").Value = 44505
$ws.Range("K3").Value = "Golden Nugget"
$ws.Range("L3").Value = "Primera"
$ws.Range("M3").Value = 50
$ws.Range("N3").Value = 15000
$ws.Range("O3").Value = 15000
$ws.Range("P3").Value = 15000
$ws.Range("Q3").Value = "$/bandeja 10 kilos"
$ws.Range("S3").Value = 1500

$ws.Range("D4").Value = 44902
$ws.Range("K4").Value = "Golden Nugget"
$ws.Range("L4").Value = "Especial"
$ws.Range("M4").Value = 60
$ws.Range("N4").Value = 15000
$ws.Range("O4").Value = 15000
$ws.Range("P4").Value = 15000
$ws.Range("Q4").Value = "$/caja 10 kilos"
$ws.Range("S4").Value = 1500

$ws.Range("D5").Value = 44902
$ws.Range("K5").Value = "Golden Nugget"
$ws.Range("L5").Value = "Primera"
$ws.Range("M5").Value = 70
$ws.Range("N5").Value = 13000
$ws.Range("O5").Value = 13000
$ws.Range("P5").Value = 13000
$ws.Range("Q5").Value = "$/caja 10 kilos"
$ws.Range("S5").Value = 1300
